$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column F with header and DB ID values
$ws.Range("F1").Value = "DB ID's"
$ws.Range("F2").Value = 3
$ws.Range("F3").Value = 11
$ws.Range("F4").Value = 5
$ws.Range("F5").Value = 6
$ws.Range("F6").Value = 12
$ws.Range("F7").Value = 4
$ws.Range("F8").Value = 7
$ws.Range("F9").Value = 8
$ws.Range("F10").Value = 9
$ws.Range("F11").Value = 10
$ws.Range("F12").Value = 14
$ws.Range("F13").Value = 15
$ws.Range("F14").Value = 16
$ws.Range("F15").Value = 17
$ws.Range("F16").Value = 18
$ws.Range("F17").Value = 19
$ws.Range("F18").Value = 20
$ws.Range("F19").Value = 21
$ws.Range("F20").Value = 22
$ws.Range("F21").Value = 23

# Update selection to H4 (matches the new active cell in the saved view)
$ws.Range("H4").Select()
